$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 110046148
$ws.Range("B8").Value = 78098
$ws.Range("E8").Value = 6453
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "Vedskivlav"
$ws.Range("F8").NumberFormat = "General"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "Hertelidea botryosa"
$ws.Range("G8").NumberFormat = "General"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("H8").NumberFormat = "General"
$ws.Range("Q8").Value = 730378.4553517678
$ws.Range("R8").Value = 7376568.376099556
$ws.Range("AX8").NumberFormat = "@"
$ws.Range("AX8").Value = "Robert Sandberg"
$ws.Range("AX8").NumberFormat = "General"
# Row 9
$ws.Range("A9").Value = 110046156
$ws.Range("B9").Value = 78098
$ws.Range("E9").Value = 6453
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "Vedskivlav"
$ws.Range("F9").NumberFormat = "General"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "Hertelidea botryosa"
$ws.Range("G9").NumberFormat = "General"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("H9").NumberFormat = "General"
$ws.Range("Q9").Value = 730373.7458103633
$ws.Range("R9").Value = 7376547.475861446
$ws.Range("AX9").NumberFormat = "@"
$ws.Range("AX9").Value = "Rasmus Häggqvist"
$ws.Range("AX9").NumberFormat = "General"
# Row 10
$ws.Range("A10").Value = 110046158
$ws.Range("B10").Value = 77258
$ws.Range("E10").Value = 6446
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "Kolflarnlav"
$ws.Range("F10").NumberFormat = "General"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "Carbonicola anthracophila"
$ws.Range("G10").NumberFormat = "General"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("H10").NumberFormat = "General"
$ws.Range("Q10").Value = 730374.7825402109
$ws.Range("R10").Value = 7376554.398445786
$ws.Range("AX10").NumberFormat = "@"
$ws.Range("AX10").Value = "Rasmus Häggqvist"
$ws.Range("AX10").NumberFormat = "General"
# Row 11
$ws.Range("A11").Value = 110046134
$ws.Range("B11").Value = 81236
$ws.Range("E11").Value = 1312
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "Gammelgransskål"
$ws.Range("F11").NumberFormat = "General"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "Pseudographis pinicola"
$ws.Range("G11").NumberFormat = "General"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "(Nyl.) Rehm"
$ws.Range("H11").NumberFormat = "General"
$ws.Range("Q11").Value = 730383.5096258806
$ws.Range("R11").Value = 7376429.645832454
$ws.Range("Y11").NumberFormat = "@"
$ws.Range("Y11").Value = "2022-06-02"
$ws.Range("Y11").NumberFormat = "General"
$ws.Range("AA11").NumberFormat = "@"
$ws.Range("AA11").Value = "2022-06-02"
$ws.Range("AA11").NumberFormat = "General"
$ws.Range("AC11").NumberFormat = "@"
$ws.Range("AC11").Value = "på flera granar"
$ws.Range("AC11").NumberFormat = "General"
$ws.Range("AX11").NumberFormat = "@"
$ws.Range("AX11").Value = "Robert Sandberg"
$ws.Range("AX11").NumberFormat = "General"
# Row 12
$ws.Range("A12").Value = 110046155
$ws.Range("B12").Value = 81236
$ws.Range("E12").Value = 1312
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "Gammelgransskål"
$ws.Range("F12").NumberFormat = "General"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "Pseudographis pinicola"
$ws.Range("G12").NumberFormat = "General"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "(Nyl.) Rehm"
$ws.Range("H12").NumberFormat = "General"
$ws.Range("Q12").Value = 730414.7718965814
$ws.Range("R12").Value = 7376382.768500465
$ws.Range("AX12").NumberFormat = "@"
$ws.Range("AX12").Value = "Rasmus Häggqvist"
$ws.Range("AX12").NumberFormat = "General"
# Row 13
$ws.Range("A13").Value = 110046052
$ws.Range("B13").Value = 77259
$ws.Range("E13").Value = 228912
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "Mörk kolflarnlav"
$ws.Range("F13").NumberFormat = "General"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "Carbonicola myrmecina"
$ws.Range("G13").NumberFormat = "General"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("H13").NumberFormat = "General"
$ws.Range("Q13").Value = 730361.5117082358
$ws.Range("R13").Value = 7376418.170894846
# Row 14
$ws.Range("A14").Value = 110046170
$ws.Range("B14").Value = 56395
$ws.Range("E14").Value = 100109
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "Tretåig hackspett"
$ws.Range("F14").NumberFormat = "General"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "Picoides tridactylus"
$ws.Range("G14").NumberFormat = "General"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "(Linnaeus, 1758)"
$ws.Range("H14").NumberFormat = "General"
$ws.Range("Q14").Value = 730426.7096293157
$ws.Range("R14").Value = 7376399.040735548
$ws.Range("AX14").NumberFormat = "@"
$ws.Range("AX14").Value = "Frédéric Forsmark"
$ws.Range("AX14").NumberFormat = "General"
# Row 15
$ws.Range("A15").Value = 110049842
$ws.Range("B15").Value = 77177
$ws.Range("E15").Value = 353
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "Dvärgbägarlav"
$ws.Range("F15").NumberFormat = "General"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "Cladonia parasitica"
$ws.Range("G15").NumberFormat = "General"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "(Hoffm.) Hoffm."
$ws.Range("H15").NumberFormat = "General"
$ws.Range("Q15").Value = 730532.734193347
$ws.Range("R15").Value = 7376460.51517869
$ws.Range("Y15").NumberFormat = "@"
$ws.Range("Y15").Value = "2022-09-15"
$ws.Range("Y15").NumberFormat = "General"
$ws.Range("AA15").NumberFormat = "@"
$ws.Range("AA15").Value = "2022-09-15"
$ws.Range("AA15").NumberFormat = "General"
$ws.Range("AX15").NumberFormat = "@"
$ws.Range("AX15").Value = "Frédéric Forsmark"
$ws.Range("AX15").NumberFormat = "General"
# Row 16
$ws.Range("A16").Value = 110046051
$ws.Range("B16").Value = 78098
$ws.Range("E16").Value = 6453
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "Vedskivlav"
$ws.Range("F16").NumberFormat = "General"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "Hertelidea botryosa"
$ws.Range("G16").NumberFormat = "General"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("H16").NumberFormat = "General"
$ws.Range("Q16").Value = 730363.3455545675
$ws.Range("R16").Value = 7376420.333717911
$ws.Range("AC16").ClearContents() | Out-Null
$ws.Range("AX16").NumberFormat = "@"
$ws.Range("AX16").Value = "Linda Spjut"
$ws.Range("AX16").NumberFormat = "General"
